$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing evaluation metric values
$ws.Range("B2").Value = 2243454.062
$ws.Range("B3").Value = 1497.816
$ws.Range("B4").Value = 1272.794

# Add new row 12 with prediction accuracy metric
$ws.Range("A12").Value = "Accuracy of correct prediction side"
$ws.Range("B12").Value = 0.451
